$wb = $excel.ActiveWorkbook
$wsRate = $wb.Worksheets.Item("rate matrix")
$wsVisual = $wb.Worksheets.Item("visual Q")

# Update the "Death" row values on both sheets
$wsRate.Range("B8").Value = 0
$wsRate.Range("H8").Value = 0

$wsVisual.Range("B8").ClearContents()
$wsVisual.Range("H8").Value = 0

# Update selections on both sheets, then make "rate matrix" the active tab
$wsVisual.Range("B9").Select()

$wsRate.Activate()
$wsRate.Range("B9").Select()
